# Config.xlsx maintenance edit:
#  - Reorder the comma-separated TikTok profile URLs stored in the
#    "ProfileUrls" row (Key/Value table on Sheet1) so @elianmita comes
#    before @selly.
#  - Move the sheet's saved active-cell selection from B10 to B8.
# Wrapped defensively so a layout change (extra/missing rows, renamed
# sheet, etc.) doesn't blow up the whole script - better try/catch logic
# and edge cases, per the usual config-sheet gremlins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    # Find the row whose Key column holds "ProfileUrls" instead of
    # hard-coding B2, so this keeps working if rows get inserted/removed.
    $profileRow = $null
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $key = $ws.Cells.Item($r, 1).Value2
        if ($key -eq "ProfileUrls") {
            $profileRow = $r
            break
        }
    }

    if ($profileRow -ne $null) {
        $valueCell = $ws.Cells.Item($profileRow, 2)
        $current = [string]$valueCell.Value2

        if ([string]::IsNullOrEmpty($current)) {
            Write-Output "ProfileUrls value is empty - nothing to reorder."
        } else {
            $urls = $current.Split(",")
            if ($urls.Count -gt 1) {
                $reversed = @()
                for ($i = $urls.Count - 1; $i -ge 0; $i--) {
                    $reversed += $urls[$i]
                }
                $valueCell.Value = [string]::Join(",", $reversed)
            } else {
                Write-Output "Only one URL present - leaving as-is."
            }
        }
    } else {
        Write-Output "Could not find a 'ProfileUrls' row - skipping URL reorder."
    }
} catch {
    Write-Output "Failed to update ProfileUrls: $_"
}

try {
    # Move the saved selection off B10 onto B8.
    $target = $ws.Range("B8")
    $target.Select()
} catch {
    Write-Output "Failed to update selection: $_"
}
